$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23; this shifts the existing rows 23-43
# down to 24-44 (preserving all their data), matching the diff which
# shows every former row's data reappearing one row lower.
$ws.Rows.Item(23).Insert()

# Populate the newly-inserted row 23 with the new data record.
$ws.Range("A23").Value = 10
$ws.Range("B23").Value = "Vega Modelo de Temuco"
$ws.Range("C23").Value = "La Araucanía"
$ws.Range("D23").Value = 44740
$ws.Range("E23").Value = 9
$ws.Range("F23").Value = 100112010
$ws.Range("G23").Value = "Achicoria"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 170
$ws.Range("K23").Value = 10000
$ws.Range("L23").Value = 11000
$ws.Range("M23").Value = 10471
$ws.Range("N23").Value = "$/caja 18 unidades"
$ws.Range("O23").Value = "Región Metropolitana"
$ws.Range("P23").Value = 582
$ws.Range("Q23").Value = 18
$ws.Range("R23").Value = "Hortaliza"
